$wb = $excel.ActiveWorkbook

# --- params sheet: num_period value 10 -> 60 ---
$wsParams = $wb.Worksheets.Item("params")
$wsParams.Range("B4").Value = 60

# --- heures_vol sheet: shift the "mois" column forward by ~4 years ---
# 2015-04 -> 2019-04 (no quote-prefix style on these rows)
$wsHeures = $wb.Worksheets.Item("heures_vol")
for ($r = 2; $r -le 61; $r++) {
    $wsHeures.Cells.Item($r, 2).Value = "2019-04"
}
# 2015-05 -> 2019-05 (quote-prefix style preserved via leading apostrophe)
for ($r = 62; $r -le 121; $r++) {
    $wsHeures.Cells.Item($r, 2).Value = "'2019-05"
}
# 2015-06 -> 2019-06
for ($r = 122; $r -le 181; $r++) {
    $wsHeures.Cells.Item($r, 2).Value = "'2019-06"
}
# 2015-07 -> 2019-07
for ($r = 182; $r -le 241; $r++) {
    $wsHeures.Cells.Item($r, 2).Value = "'2019-07"
}

# --- selections / active sheet bookkeeping ---
# avions was the active sheet/tab (tabSelected) with selection C6 (kept as-is).
$wsAvions = $wb.Worksheets.Item("avions")
$null = $wsAvions.Activate()
$null = $wsAvions.Range("C6").Select()

# heures_vol gains a selection at B2 (no longer the active tab afterwards).
$null = $wsHeures.Activate()
$null = $wsHeures.Range("B2").Select()

# params becomes the active tab, with selection moved to B5.
$null = $wsParams.Activate()
$null = $wsParams.Range("B5").Select()
